$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the current "Address" column (column C)
$ws.Range("C1:E1").EntireColumn.Insert()

# Rename the first two (now still A/B) headers
$ws.Range("A1").Value = "Offer Quantity *"
$ws.Range("B1").Value = "User (email) *"

# Fill in the new header cells
$ws.Range("C1").Value = "First Name *"
$ws.Range("D1").Value = "Middle Name"
$ws.Range("E1").Value = "Last Name *"

# Fill in First / Last name sample data (Middle Name stays blank)
$ws.Range("C2").Value = "Emp1"
$ws.Range("E2").Value = "John"

$ws.Range("C3").Value = "Emp2"
$ws.Range("E3").Value = "James"

$ws.Range("C4").Value = "Emp3"
$ws.Range("E4").Value = "Jim"

$ws.Range("C5").Value = "Emp4"
$ws.Range("E5").Value = "Jack"

# View changes: scroll to E1, select K1
$ws.Range("K1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5
